# Append a new price-tracking row (2026-02-07) to the bottom of the sheet,
# mirroring the existing Date/Price/Discount/Incredible rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the existing data (row 38 here).
$newRow = $ws.Cells.Item($ws.UsedRange.Rows.Count + 1, 1).Row

$newDate = "2026-02-07"
$newPrice = "324800"
$newDiscount = "44"
$newIncredible = "0"

$rowRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 4))

# Force the cells to be treated as text so values like "324800" / "44" are
# stored as plain strings (matching the rest of the column) instead of being
# auto-converted to numbers/dates by Excel.
$rowRange.NumberFormat = "@"

$ws.Cells.Item($newRow, 1).Value = $newDate
$ws.Cells.Item($newRow, 2).Value = $newPrice
$ws.Cells.Item($newRow, 3).Value = $newDiscount
$ws.Cells.Item($newRow, 4).Value = $newIncredible

# Drop the temporary text formatting again so the new cells fall back to the
# workbook's default (unstyled) cell formatting, just like the pre-existing
# rows.
$rowRange.ClearFormats()
